$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(3).TextFrame.TextRange

# The paragraph currently holds the sentence split across several runs
# ("Followed" / " " / "by" / " " / "a" / " " / "picture"). Their
# concatenation already equals the target string, so assigning that same
# value directly is a no-op for the host's text-diffing and the runs
# would stay split. Route through an unrelated placeholder value first
# (sharing no prefix/suffix with the target, so the host can't patch the
# existing runs in place) which forces the paragraph to collapse down to
# a single run; then set the final merged text onto that run.
$tr.Text = "zzzzzzzzzzzzzzzzzzzz"
$tr.Text = "Followed by a picture"
